$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Friday hours for the week of row 7 (3.25 -> 6.25)
$ws.Range("H7").Value = 6.25

# Update the active cell selection to O11 (as recorded in the sheet view)
$ws.Range("O11").Select()
